# Auto-generated Excel COM-interop script
# Applies the Moogle_Profits profit-recalculation data refresh
# (chore: update Sheets via scheduled runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 32940.7
$ws.Cells.Item(74, 9).Value = 32940.7
$ws.Cells.Item(74, 11).Value = 32940.7
$ws.Cells.Item(74, 13).Value = -32004.7
$ws.Cells.Item(77, 8).Value = 32940.7
$ws.Cells.Item(77, 9).Value = 32940.7
$ws.Cells.Item(77, 11).Value = 164703.5
$ws.Cells.Item(77, 13).Value = -160023.5
$ws.Cells.Item(92, 8).Value = 1556.9
$ws.Cells.Item(92, 9).Value = 1227.4375
$ws.Cells.Item(92, 10).Value = 2874.75
$ws.Cells.Item(92, 11).Value = 1227.4375
$ws.Cells.Item(92, 12).Value = 2874.75
$ws.Cells.Item(92, 13).Value = 20.5625
$ws.Cells.Item(92, 14).Value = -5370.75
$ws.Cells.Item(125, 8).Value = 170227.67
$ws.Cells.Item(125, 9).Value = 4082.5
$ws.Cells.Item(125, 10).Value = 502518
$ws.Cells.Item(125, 11).Value = 36742.5
$ws.Cells.Item(125, 12).Value = 4522662
$ws.Cells.Item(125, 13).Value = -34282.5
$ws.Cells.Item(125, 14).Value = -4527582
$ws.Cells.Item(132, 8).Value = 1743.5227
$ws.Cells.Item(132, 9).Value = 1743.5227
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 5230.5681
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -2700.5681
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 2179.9033
$ws.Cells.Item(138, 9).Value = 1862.9259
$ws.Cells.Item(138, 11).Value = 5588.7777
$ws.Cells.Item(138, 13).Value = -448.7776999999996
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6552.207
$ws.Cells.Item(32, 9).Value = 3154.6924
$ws.Cells.Item(32, 11).Value = 3154.6924
$ws.Cells.Item(32, 13).Value = -2867.6924
$ws.Cells.Item(45, 8).Value = 2546.5186
$ws.Cells.Item(45, 9).Value = 2243.6365
$ws.Cells.Item(45, 10).Value = 3879.2
$ws.Cells.Item(45, 11).Value = 2243.6365
$ws.Cells.Item(45, 12).Value = 3879.2
$ws.Cells.Item(45, 13).Value = -1866.6365
$ws.Cells.Item(45, 14).Value = -4633.2
$ws.Cells.Item(61, 8).Value = 4908.1514
$ws.Cells.Item(61, 10).Value = 5485.4736
$ws.Cells.Item(61, 12).Value = 5485.4736
$ws.Cells.Item(61, 14).Value = -5909.4736
$ws.Cells.Item(88, 8).Value = 3670.8
$ws.Cells.Item(88, 9).Value = 4334.6665
$ws.Cells.Item(88, 10).Value = 2675
$ws.Cells.Item(88, 11).Value = 4334.6665
$ws.Cells.Item(88, 12).Value = 2675
$ws.Cells.Item(88, 13).Value = -3928.6665
$ws.Cells.Item(88, 14).Value = -3487
$ws.Cells.Item(91, 8).Value = 3670.8
$ws.Cells.Item(91, 9).Value = 4334.6665
$ws.Cells.Item(91, 10).Value = 2675
$ws.Cells.Item(91, 11).Value = 4334.6665
$ws.Cells.Item(91, 12).Value = 2675
$ws.Cells.Item(91, 13).Value = -2930.6665
$ws.Cells.Item(91, 14).Value = -5483
$ws.Cells.Item(132, 8).Value = 2102.3713
$ws.Cells.Item(132, 9).Value = 1424.3448
$ws.Cells.Item(132, 11).Value = 4273.0344
$ws.Cells.Item(132, 13).Value = -1743.0344
$ws.Cells.Item(136, 8).Value = 4908.1514
$ws.Cells.Item(136, 10).Value = 5485.4736
$ws.Cells.Item(136, 12).Value = 16456.4208
$ws.Cells.Item(136, 14).Value = -21556.4208
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1787.8125
$ws.Cells.Item(86, 9).Value = 1588.5555
$ws.Cells.Item(86, 10).Value = 2044
$ws.Cells.Item(86, 11).Value = 1588.5555
$ws.Cells.Item(86, 12).Value = 2044
$ws.Cells.Item(86, 13).Value = -465.5554999999999
$ws.Cells.Item(86, 14).Value = -4290
$ws.Cells.Item(89, 8).Value = 1787.8125
$ws.Cells.Item(89, 9).Value = 1588.5555
$ws.Cells.Item(89, 10).Value = 2044
$ws.Cells.Item(89, 11).Value = 7942.7775
$ws.Cells.Item(89, 12).Value = 10220
$ws.Cells.Item(89, 13).Value = -2326.7775
$ws.Cells.Item(89, 14).Value = -21452
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6464.4116
$ws.Cells.Item(31, 9).Value = 4010.95
$ws.Cells.Item(31, 10).Value = 8047.2905
$ws.Cells.Item(31, 11).Value = 4010.95
$ws.Cells.Item(31, 12).Value = 8047.2905
$ws.Cells.Item(31, 13).Value = -3715.95
$ws.Cells.Item(31, 14).Value = -8637.290499999999
$ws.Cells.Item(34, 8).Value = 6464.4116
$ws.Cells.Item(34, 9).Value = 4010.95
$ws.Cells.Item(34, 10).Value = 8047.2905
$ws.Cells.Item(34, 11).Value = 4010.95
$ws.Cells.Item(34, 12).Value = 8047.2905
$ws.Cells.Item(34, 13).Value = -3808.95
$ws.Cells.Item(34, 14).Value = -8451.290499999999
$ws.Cells.Item(62, 8).Value = 6395.579
$ws.Cells.Item(62, 9).Value = 4702.5713
$ws.Cells.Item(62, 10).Value = 7383.1665
$ws.Cells.Item(62, 11).Value = 4702.5713
$ws.Cells.Item(62, 12).Value = 7383.1665
$ws.Cells.Item(62, 13).Value = -4078.5713
$ws.Cells.Item(62, 14).Value = -8631.166499999999
$ws.Cells.Item(65, 8).Value = 6395.579
$ws.Cells.Item(65, 9).Value = 4702.5713
$ws.Cells.Item(65, 10).Value = 7383.1665
$ws.Cells.Item(65, 11).Value = 23512.8565
$ws.Cells.Item(65, 12).Value = 36915.8325
$ws.Cells.Item(65, 13).Value = -20392.8565
$ws.Cells.Item(65, 14).Value = -43155.8325
$ws.Cells.Item(132, 8).Value = 2430.5293
$ws.Cells.Item(132, 9).Value = 1412.8837
$ws.Cells.Item(132, 11).Value = 4238.6511
$ws.Cells.Item(132, 13).Value = -1708.6511
$ws.Cells.Item(134, 8).Value = 1962.037
$ws.Cells.Item(134, 9).Value = 1682.875
$ws.Cells.Item(134, 11).Value = 5048.625
$ws.Cells.Item(134, 13).Value = -2513.625
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 11499.25
$ws.Cells.Item(55, 10).Value = 11499.25
$ws.Cells.Item(55, 12).Value = 34497.75
$ws.Cells.Item(55, 14).Value = -34851.75
$ws.Cells.Item(121, 8).Value = 906063.0600000001
$ws.Cells.Item(121, 10).Value = 1266924.9
$ws.Cells.Item(121, 12).Value = 3800774.7
$ws.Cells.Item(121, 14).Value = -3803394.7
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 67499.5
$ws.Cells.Item(15, 10).Value = 67499.5
$ws.Cells.Item(15, 12).Value = 67499.5
$ws.Cells.Item(15, 14).Value = -68075.5
$ws.Cells.Item(81, 8).Value = 67499.5
$ws.Cells.Item(81, 10).Value = 67499.5
$ws.Cells.Item(81, 12).Value = 67499.5
$ws.Cells.Item(81, 14).Value = -69495.5
$ws.Cells.Item(84, 8).Value = 67499.5
$ws.Cells.Item(84, 10).Value = 67499.5
$ws.Cells.Item(84, 12).Value = 202498.5
$ws.Cells.Item(84, 14).Value = -212482.5
$ws.Cells.Item(120, 8).Value = 58998
$ws.Cells.Item(120, 10).Value = 58998
$ws.Cells.Item(120, 12).Value = 58998
$ws.Cells.Item(120, 14).Value = -68674
$ws.Cells.Item(121, 8).Value = 49000
$ws.Cells.Item(121, 10).Value = 49000
$ws.Cells.Item(121, 12).Value = 49000
$ws.Cells.Item(121, 14).Value = -52494
$ws.Cells.Item(132, 8).Value = 4152.851
$ws.Cells.Item(132, 9).Value = 2568.4333
$ws.Cells.Item(132, 11).Value = 7705.2999
$ws.Cells.Item(132, 13).Value = -5175.2999
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 11351
$ws.Cells.Item(7, 9).Value = 6801.3335
$ws.Cells.Item(7, 11).Value = 6801.3335
$ws.Cells.Item(7, 13).Value = -6689.3335
$ws.Cells.Item(68, 8).Value = 6611.773
$ws.Cells.Item(68, 9).Value = 4437.769
$ws.Cells.Item(68, 11).Value = 4437.769
$ws.Cells.Item(68, 13).Value = -3688.769
$ws.Cells.Item(71, 8).Value = 6611.773
$ws.Cells.Item(71, 9).Value = 4437.769
$ws.Cells.Item(71, 11).Value = 22188.845
$ws.Cells.Item(71, 13).Value = -18444.845
$ws.Cells.Item(82, 8).Value = 1757.421
$ws.Cells.Item(82, 9).Value = 945.1
$ws.Cells.Item(82, 10).Value = 2660
$ws.Cells.Item(82, 11).Value = 945.1
$ws.Cells.Item(82, 12).Value = 2660
$ws.Cells.Item(82, 13).Value = -584.1
$ws.Cells.Item(82, 14).Value = -3382
$ws.Cells.Item(85, 8).Value = 1757.421
$ws.Cells.Item(85, 9).Value = 945.1
$ws.Cells.Item(85, 10).Value = 2660
$ws.Cells.Item(85, 11).Value = 945.1
$ws.Cells.Item(85, 12).Value = 2660
$ws.Cells.Item(85, 13).Value = 302.9
$ws.Cells.Item(85, 14).Value = -5156
$ws.Cells.Item(93, 8).Value = 3294.1333
$ws.Cells.Item(93, 9).Value = 2232.6667
$ws.Cells.Item(93, 11).Value = 2232.6667
$ws.Cells.Item(93, 13).Value = -984.6667000000002
$ws.Cells.Item(94, 8).Value = 90737.25
$ws.Cells.Item(94, 10).Value = 90737.25
$ws.Cells.Item(94, 12).Value = 90737.25
$ws.Cells.Item(94, 14).Value = -92089.25
$ws.Cells.Item(126, 8).Value = 11351
$ws.Cells.Item(126, 9).Value = 6801.3335
$ws.Cells.Item(126, 11).Value = 20404.0005
$ws.Cells.Item(126, 13).Value = -17934.0005
$ws.Cells.Item(132, 8).Value = 6688.5747
$ws.Cells.Item(132, 9).Value = 5349.893
$ws.Cells.Item(132, 10).Value = 8661.368
$ws.Cells.Item(132, 11).Value = 16049.679
$ws.Cells.Item(132, 12).Value = 25984.104
$ws.Cells.Item(132, 13).Value = -13519.679
$ws.Cells.Item(132, 14).Value = -31044.104
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 11100
$ws.Cells.Item(96, 10).Value = 15200
$ws.Cells.Item(96, 12).Value = 15200
$ws.Cells.Item(96, 14).Value = -17946
$ws.Cells.Item(124, 8).Value = 53749.75
$ws.Cells.Item(124, 10).Value = 53749.75
$ws.Cells.Item(124, 12).Value = 53749.75
$ws.Cells.Item(124, 14).Value = -63569.75
$ws.Cells.Item(132, 8).Value = 2236.1277
$ws.Cells.Item(132, 9).Value = 1862.6216
$ws.Cells.Item(132, 11).Value = 5587.864799999999
$ws.Cells.Item(132, 13).Value = -3057.864799999999
$ws.Cells.Item(136, 8).Value = 4184.433
$ws.Cells.Item(136, 9).Value = 4190.45
$ws.Cells.Item(136, 10).Value = 4172.4
$ws.Cells.Item(136, 11).Value = 12571.35
$ws.Cells.Item(136, 12).Value = 12517.2
$ws.Cells.Item(136, 13).Value = -10021.35
$ws.Cells.Item(136, 14).Value = -17617.2
